# Auto-generated Excel COM-interop script applying Goblin_Profits market-data refresh
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(6, 8).Value = 100000160
$ws.Cells.Item(6, 9).Value = 100000160
$ws.Cells.Item(6, 10).Value = 0
$ws.Cells.Item(6, 11).Value = 300000480
$ws.Cells.Item(6, 12).Value = 0
$ws.Cells.Item(6, 13).Value = -300000368
$ws.Cells.Item(6, 14).ClearContents()
$ws.Cells.Item(29, 8).Value = 166668020
$ws.Cells.Item(29, 10).Value = 4003
$ws.Cells.Item(29, 12).Value = 12009
$ws.Cells.Item(29, 14).Value = -12571
$ws.Cells.Item(38, 8).Value = 1647.1538
$ws.Cells.Item(38, 9).Value = 601.8889
$ws.Cells.Item(38, 10).Value = 3999
$ws.Cells.Item(38, 11).Value = 1805.6667
$ws.Cells.Item(38, 12).Value = 11997
$ws.Cells.Item(38, 13).Value = -1433.6667
$ws.Cells.Item(38, 14).Value = -12741
$ws.Cells.Item(58, 8).Value = 35715572
$ws.Cells.Item(58, 9).Value = 41667170
$ws.Cells.Item(58, 10).Value = 6000
$ws.Cells.Item(58, 11).Value = 125001510
$ws.Cells.Item(58, 12).Value = 18000
$ws.Cells.Item(58, 13).Value = -125001360
$ws.Cells.Item(58, 14).Value = -18300
$ws.Cells.Item(80, 8).Value = 1114.3334
$ws.Cells.Item(80, 10).Value = 1153.75
$ws.Cells.Item(80, 12).Value = 3461.25
$ws.Cells.Item(80, 14).Value = -5457.25
$ws.Cells.Item(83, 8).Value = 1114.3334
$ws.Cells.Item(83, 10).Value = 1153.75
$ws.Cells.Item(83, 12).Value = 10383.75
$ws.Cells.Item(83, 14).Value = -20367.75
$ws.Cells.Item(92, 8).Value = 1524.1305
$ws.Cells.Item(92, 9).Value = 1528.6471
$ws.Cells.Item(92, 11).Value = 1528.6471
$ws.Cells.Item(92, 13).Value = -280.6470999999999
$ws.Cells.Item(98, 8).Value = 9696.157999999999
$ws.Cells.Item(98, 9).Value = 9696.157999999999
$ws.Cells.Item(98, 11).Value = 9696.157999999999
$ws.Cells.Item(98, 13).Value = -8198.157999999999
$ws.Cells.Item(122, 8).Value = 9696.157999999999
$ws.Cells.Item(122, 9).Value = 9696.157999999999
$ws.Cells.Item(122, 11).Value = 29088.474
$ws.Cells.Item(122, 13).Value = -26638.474

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 3053.8
$ws.Cells.Item(32, 9).Value = 2737.6191
$ws.Cells.Item(32, 10).Value = 7480.3335
$ws.Cells.Item(32, 11).Value = 2737.6191
$ws.Cells.Item(32, 12).Value = 7480.3335
$ws.Cells.Item(32, 13).Value = -2450.6191
$ws.Cells.Item(32, 14).Value = -8054.3335
$ws.Cells.Item(45, 8).Value = 2284.3333
$ws.Cells.Item(45, 9).Value = 2037.4546
$ws.Cells.Item(45, 11).Value = 2037.4546
$ws.Cells.Item(45, 13).Value = -1660.4546
$ws.Cells.Item(74, 8).Value = 1866.5
$ws.Cells.Item(74, 9).Value = 1490.7727
$ws.Cells.Item(74, 11).Value = 1490.7727
$ws.Cells.Item(74, 13).Value = -616.7727
$ws.Cells.Item(77, 8).Value = 1866.5
$ws.Cells.Item(77, 9).Value = 1490.7727
$ws.Cells.Item(77, 11).Value = 7453.863499999999
$ws.Cells.Item(77, 13).Value = -3085.863499999999
$ws.Cells.Item(102, 8).Value = 4489.2666
$ws.Cells.Item(102, 9).Value = 2413.476
$ws.Cells.Item(102, 11).Value = 2413.476
$ws.Cells.Item(102, 13).Value = -791.4760000000001
$ws.Cells.Item(110, 8).Value = 1612.8462
$ws.Cells.Item(110, 9).Value = 706.1111
$ws.Cells.Item(110, 11).Value = 706.1111
$ws.Cells.Item(110, 13).Value = 1338.8889
$ws.Cells.Item(122, 8).Value = 2380.1667
$ws.Cells.Item(122, 9).Value = 2173.5557
$ws.Cells.Item(122, 11).Value = 6520.6671
$ws.Cells.Item(122, 13).Value = -4070.6671
$ws.Cells.Item(132, 8).Value = 1231.7435
$ws.Cells.Item(132, 9).Value = 1231.7435
$ws.Cells.Item(132, 11).Value = 3695.2305
$ws.Cells.Item(132, 13).Value = -1165.2305

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 1155.4642
$ws.Cells.Item(20, 9).Value = 833.17645
$ws.Cells.Item(20, 11).Value = 833.17645
$ws.Cells.Item(20, 13).Value = -586.17645
$ws.Cells.Item(86, 8).Value = 560391.9
$ws.Cells.Item(86, 9).Value = 1574.4
$ws.Cells.Item(86, 10).Value = 2112662.8
$ws.Cells.Item(86, 11).Value = 1574.4
$ws.Cells.Item(86, 12).Value = 2112662.8
$ws.Cells.Item(86, 13).Value = -451.4000000000001
$ws.Cells.Item(86, 14).Value = -2114908.8
$ws.Cells.Item(89, 8).Value = 560391.9
$ws.Cells.Item(89, 9).Value = 1574.4
$ws.Cells.Item(89, 10).Value = 2112662.8
$ws.Cells.Item(89, 11).Value = 7872
$ws.Cells.Item(89, 12).Value = 10563314
$ws.Cells.Item(89, 13).Value = -2256
$ws.Cells.Item(89, 14).Value = -10574546
$ws.Cells.Item(94, 8).Value = 982.17645
$ws.Cells.Item(94, 9).Value = 1105.1538
$ws.Cells.Item(94, 10).Value = 582.5
$ws.Cells.Item(94, 11).Value = 1105.1538
$ws.Cells.Item(94, 12).Value = 582.5
$ws.Cells.Item(94, 13).Value = -654.1538
$ws.Cells.Item(94, 14).Value = -1484.5
$ws.Cells.Item(105, 8).Value = 2521.4119
$ws.Cells.Item(105, 9).Value = 2366.5625
$ws.Cells.Item(105, 11).Value = 2366.5625
$ws.Cells.Item(105, 13).Value = -619.5625
$ws.Cells.Item(107, 8).Value = 3738.0667
$ws.Cells.Item(107, 10).Value = 6370
$ws.Cells.Item(107, 12).Value = 6370
$ws.Cells.Item(107, 14).Value = -10210
$ws.Cells.Item(134, 8).Value = 2510.7
$ws.Cells.Item(134, 9).Value = 2379.2173
$ws.Cells.Item(134, 11).Value = 7137.651899999999
$ws.Cells.Item(134, 13).Value = -4602.651899999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 3701.6553
$ws.Cells.Item(31, 9).Value = 1585.125
$ws.Cells.Item(31, 11).Value = 1585.125
$ws.Cells.Item(31, 13).Value = -1290.125
$ws.Cells.Item(34, 8).Value = 3701.6553
$ws.Cells.Item(34, 9).Value = 1585.125
$ws.Cells.Item(34, 11).Value = 1585.125
$ws.Cells.Item(34, 13).Value = -1383.125
$ws.Cells.Item(56, 8).Value = 8597.5
$ws.Cells.Item(56, 9).Value = 8597.5
$ws.Cells.Item(56, 11).Value = 8597.5
$ws.Cells.Item(56, 13).Value = -7752.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(70, 8).Value = 552.5
$ws.Cells.Item(70, 9).Value = 552.5
$ws.Cells.Item(70, 11).Value = 1657.5
$ws.Cells.Item(70, 13).Value = -1342.5
$ws.Cells.Item(73, 8).Value = 552.5
$ws.Cells.Item(73, 9).Value = 552.5
$ws.Cells.Item(73, 11).Value = 1657.5
$ws.Cells.Item(73, 13).Value = -565.5
$ws.Cells.Item(113, 8).Value = 1588.9131
$ws.Cells.Item(113, 9).Value = 598.6667
$ws.Cells.Item(113, 10).Value = 1938.4117
$ws.Cells.Item(113, 11).Value = 1796.0001
$ws.Cells.Item(113, 12).Value = 5815.2351
$ws.Cells.Item(113, 13).Value = 373.9999
$ws.Cells.Item(113, 14).Value = -10155.2351
$ws.Cells.Item(129, 8).Value = 3416
$ws.Cells.Item(129, 9).Value = 336.2
$ws.Cells.Item(129, 11).Value = 1008.6
$ws.Cells.Item(129, 13).Value = 3991.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 6718.375
$ws.Cells.Item(70, 9).Value = 6707.8335
$ws.Cells.Item(70, 10).Value = 6750
$ws.Cells.Item(70, 11).Value = 6707.8335
$ws.Cells.Item(70, 12).Value = 6750
$ws.Cells.Item(70, 13).Value = -6437.8335
$ws.Cells.Item(70, 14).Value = -7290
$ws.Cells.Item(73, 8).Value = 6718.375
$ws.Cells.Item(73, 9).Value = 6707.8335
$ws.Cells.Item(73, 10).Value = 6750
$ws.Cells.Item(73, 11).Value = 6707.8335
$ws.Cells.Item(73, 12).Value = 6750
$ws.Cells.Item(73, 13).Value = -5771.8335
$ws.Cells.Item(73, 14).Value = -8622
$ws.Cells.Item(102, 8).Value = 5716.923
$ws.Cells.Item(102, 9).Value = 4258.8335
$ws.Cells.Item(102, 11).Value = 4258.8335
$ws.Cells.Item(102, 13).Value = -2636.8335
$ws.Cells.Item(113, 8).Value = 7504.9565
$ws.Cells.Item(113, 9).Value = 4662.9165
$ws.Cells.Item(113, 11).Value = 4662.9165
$ws.Cells.Item(113, 13).Value = -2492.9165
$ws.Cells.Item(122, 8).Value = 2870.0908
$ws.Cells.Item(122, 9).Value = 2191.8
$ws.Cells.Item(122, 10).Value = 3435.3333
$ws.Cells.Item(122, 11).Value = 6575.400000000001
$ws.Cells.Item(122, 12).Value = 10305.9999
$ws.Cells.Item(122, 13).Value = -4125.400000000001
$ws.Cells.Item(122, 14).Value = -15205.9999
$ws.Cells.Item(132, 8).Value = 2314.0588
$ws.Cells.Item(132, 9).Value = 1966.76
$ws.Cells.Item(132, 10).Value = 3278.7778
$ws.Cells.Item(132, 11).Value = 5900.28
$ws.Cells.Item(132, 12).Value = 9836.3334
$ws.Cells.Item(132, 13).Value = -3370.28
$ws.Cells.Item(132, 14).Value = -14896.3334

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(82, 8).Value = 2681.0588
$ws.Cells.Item(82, 10).Value = 4597.75
$ws.Cells.Item(82, 12).Value = 4597.75
$ws.Cells.Item(82, 14).Value = -5319.75
$ws.Cells.Item(85, 8).Value = 2681.0588
$ws.Cells.Item(85, 10).Value = 4597.75
$ws.Cells.Item(85, 12).Value = 4597.75
$ws.Cells.Item(85, 14).Value = -7093.75
$ws.Cells.Item(94, 8).Value = 32999.5
$ws.Cells.Item(94, 10).Value = 32999.5
$ws.Cells.Item(94, 12).Value = 32999.5
$ws.Cells.Item(94, 14).Value = -34351.5
$ws.Cells.Item(132, 8).Value = 3000.9473
$ws.Cells.Item(132, 9).Value = 2834.3333
$ws.Cells.Item(132, 11).Value = 8502.999899999999
$ws.Cells.Item(132, 13).Value = -5972.999899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(96, 8).Value = 4428.4287
$ws.Cells.Item(96, 9).Value = 3500
$ws.Cells.Item(96, 10).Value = 4799.8
$ws.Cells.Item(96, 11).Value = 3500
$ws.Cells.Item(96, 12).Value = 4799.8
$ws.Cells.Item(96, 13).Value = -2127
$ws.Cells.Item(96, 14).Value = -7545.8
$ws.Cells.Item(126, 8).Value = 3327.8125
$ws.Cells.Item(126, 9).Value = 2857.5
$ws.Cells.Item(126, 10).Value = 4738.75
$ws.Cells.Item(126, 11).Value = 8572.5
$ws.Cells.Item(126, 12).Value = 14216.25
$ws.Cells.Item(126, 13).Value = -6102.5
$ws.Cells.Item(126, 14).Value = -19156.25
$ws.Cells.Item(132, 8).Value = 3799.2273
$ws.Cells.Item(132, 9).Value = 3505.9443
$ws.Cells.Item(132, 11).Value = 10517.8329
$ws.Cells.Item(132, 13).Value = -7987.832900000001
